$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Acquisition Date value that was in I16 (keep cell present, but empty)
$ws.Range("I16").Value = ""

# Normalize the scientific-notation "Current Balance" text values in column M
$ws.Range("M27").Value = "9.40E+13"
$ws.Range("M28").Value = "9.40E+13"
$ws.Range("M29").Value = "9.40E+13"
$ws.Range("M31").Value = "9.40E+13"
$ws.Range("M32").Value = "9.40E+13"
$ws.Range("M33").Value = "9.40E+13"
$ws.Range("M34").Value = "9.41E+13"

# Remove the TOTAL row (row 35) entirely
$ws.Rows.Item(35).Delete()
